$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.174.65"
$ws.Range("E2").Value = "  +4.30%  "
$ws.Range("D3").Value = "'1.907.75"
$ws.Range("E3").Value = "  +5.25%  "
$ws.Range("D4").Value = "'0.9986"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'252.79"
$ws.Range("E5").Value = "  +2.69%  "
$ws.Range("D6").Value = "'0.9985"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "'0.5096"
$ws.Range("E7").Value = "  +3.29%  "
$ws.Range("D8").Value = "'45.19"
$ws.Range("E8").Value = "  +3.62%  "
$ws.Range("D9").Value = "'0.3023"
$ws.Range("E9").Value = "  +8.98%  "
$ws.Range("D10").Value = "'0.06822"
$ws.Range("E10").Value = "  +6.57%  "
$ws.Range("D11").Value = "'1.905.58"
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("D13").Value = "'0.07323"
$ws.Range("E13").Value = "  +3.54%  "
$ws.Range("E14").Value = "  +8.11%  "
$ws.Range("D15").Value = "'87.16"
$ws.Range("E15").Value = "  +4.08%  "
$ws.Range("D16").Value = "'4.927"
$ws.Range("E16").Value = "  +5.38%  "
$ws.Range("D17").Value = "'30.168.91"
$ws.Range("E17").Value = "  +4.15%  "
$ws.Range("D18").Value = "'0.000008296"
$ws.Range("E18").Value = "  +13.28%  "
$ws.Range("D19").Value = "'0.9983"
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").Value = "'13.10"
$ws.Range("E20").Value = "  +6.96%  "
$ws.Range("D21").Value = "'2.150.85"
$ws.Range("E21").Value = "  +5.34%  "
$ws.Range("D22").Value = "'0.9980"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'4.831"
$ws.Range("E23").Value = "  +5.68%  "
$ws.Range("D24").Value = "'5.770"
$ws.Range("E24").Value = "  +8.04%  "
$ws.Range("D25").Value = "'9.365"
$ws.Range("E25").Value = "  +6.41%  "
$ws.Range("D26").Value = "'148.00"
$ws.Range("E26").Value = "  +3.43%  "
$ws.Range("D27").Value = "'135.05"
$ws.Range("E27").Value = "  +4.53%  "
$ws.Range("D28").Value = "'17.19"
$ws.Range("E28").Value = "  +4.64%  "
$ws.Range("D29").Value = "'2.009"
$ws.Range("E29").Value = "  +6.48%  "
$ws.Range("E30").Value = "  -0.73%  "
$ws.Range("D31").Value = "'4.312"
$ws.Range("D32").Value = "'0.08899"
$ws.Range("E32").Value = "  +6.62%  "
$ws.Range("D33").Value = "'4.016"
$ws.Range("E33").Value = "  +6.40%  "
$ws.Range("D34").Value = "'0.05079"
$ws.Range("E34").Value = "  +2.78%  "
$ws.Range("E35").Value = "  +4.89%  "
$ws.Range("D36").Value = "'0.7233"
$ws.Range("E36").Value = "  +7.80%  "
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'2.305"
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("E39").Value = "  +2.76%  "
$ws.Range("D40").Value = "'0.9622"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").Value = "'0.01692"
$ws.Range("E41").Value = "  +6.85%  "
$ws.Range("D42").Value = "'6.074"
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").Value = "'0.4330"
$ws.Range("E43").Value = "  +6.23%  "
$ws.Range("D44").Value = "'105.38"
$ws.Range("E44").Value = "  +4.70%  "
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("D46").Value = "'7.708"
$ws.Range("E46").Value = "  +8.26%  "
$ws.Range("E47").Value = "  +5.18%  "
$ws.Range("E48").Value = "  +4.23%  "
$ws.Range("D49").Value = "'33.55"
$ws.Range("E49").Value = "  +6.05%  "
$ws.Range("D50").Value = "'8.474"
$ws.Range("E50").Value = "  +4.16%  "
$ws.Range("D51").Value = "'0.3830"
$ws.Range("E51").Value = "  +6.00%  "
